$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(6, 1).Value = 2
$ws.Cells.Item(7, 1).Value = 1
$ws.Cells.Item(14, 1).Value = 2
$ws.Cells.Item(15, 1).Value = 2
$ws.Cells.Item(17, 1).Value = 1
$ws.Cells.Item(19, 1).Value = 2
$ws.Cells.Item(20, 1).Value = 2
$ws.Cells.Item(21, 1).Value = 1
$ws.Cells.Item(22, 1).Value = 1
$ws.Cells.Item(24, 1).Value = 2
$ws.Cells.Item(27, 1).Value = 1
$ws.Cells.Item(28, 1).Value = 2
$ws.Cells.Item(32, 1).Value = 1
$ws.Cells.Item(34, 1).Value = 2
$ws.Cells.Item(35, 1).Value = 2
$ws.Cells.Item(36, 1).Value = 2
$ws.Cells.Item(37, 1).Value = 2
$ws.Cells.Item(38, 1).Value = 1
$ws.Cells.Item(40, 1).Value = 1
$ws.Cells.Item(41, 1).Value = 2
$ws.Cells.Item(42, 1).Value = 2
$ws.Cells.Item(44, 1).Value = 1
$ws.Cells.Item(45, 1).Value = 1
$ws.Cells.Item(49, 1).Value = 2
$ws.Cells.Item(50, 1).Value = 1
$ws.Cells.Item(51, 1).Value = 2
$ws.Cells.Item(54, 1).Value = 2
$ws.Cells.Item(55, 1).Value = 1
$ws.Cells.Item(56, 1).Value = 1
$ws.Cells.Item(58, 1).Value = 1
$ws.Cells.Item(60, 1).Value = 1
$ws.Cells.Item(62, 1).Value = 1
$ws.Cells.Item(63, 1).Value = 1
$ws.Cells.Item(65, 1).Value = 2
$ws.Cells.Item(67, 1).Value = 2
$ws.Cells.Item(68, 1).Value = 1
$ws.Cells.Item(71, 1).Value = 1
$ws.Cells.Item(76, 1).Value = 1
$ws.Cells.Item(77, 1).Value = 1
$ws.Cells.Item(78, 1).Value = 1
$ws.Cells.Item(83, 1).Value = 2
$ws.Cells.Item(84, 1).Value = 2
$ws.Cells.Item(85, 1).Value = 1
$ws.Cells.Item(86, 1).Value = 2
$ws.Cells.Item(88, 1).Value = 1
$ws.Cells.Item(91, 1).Value = 2
$ws.Cells.Item(92, 1).Value = 1
$ws.Cells.Item(94, 1).Value = 1
$ws.Cells.Item(95, 1).Value = 1
$ws.Cells.Item(104, 1).Value = 1
$ws.Cells.Item(105, 1).Value = 2
$ws.Cells.Item(107, 1).Value = 2
$ws.Cells.Item(108, 1).Value = 1
$ws.Cells.Item(109, 1).Value = 2
$ws.Cells.Item(111, 1).Value = 1
$ws.Cells.Item(116, 1).Value = 1
$ws.Cells.Item(119, 1).Value = 1
$ws.Cells.Item(124, 1).Value = 2
$ws.Cells.Item(126, 1).Value = 1
$ws.Cells.Item(127, 1).Value = 1
$ws.Cells.Item(128, 1).Value = 2
$ws.Cells.Item(129, 1).Value = 1
$ws.Cells.Item(131, 1).Value = 1
$ws.Cells.Item(134, 1).Value = 1
$ws.Cells.Item(136, 1).Value = 1
$ws.Cells.Item(140, 1).Value = 2
$ws.Cells.Item(145, 1).Value = 1
$ws.Cells.Item(146, 1).Value = 2
$ws.Cells.Item(147, 1).Value = 1
$ws.Cells.Item(148, 1).Value = 1
$ws.Cells.Item(151, 1).Value = 2
$ws.Cells.Item(152, 1).Value = 1
$ws.Cells.Item(154, 1).Value = 1
$ws.Cells.Item(155, 1).Value = 2
$ws.Cells.Item(156, 1).Value = 1
$ws.Cells.Item(157, 1).Value = 2
$ws.Cells.Item(159, 1).Value = 1
$ws.Cells.Item(161, 1).Value = 1
$ws.Cells.Item(163, 1).Value = 1
$ws.Cells.Item(169, 1).Value = 1
$ws.Cells.Item(171, 1).Value = 1
$ws.Cells.Item(177, 1).Value = 2
$ws.Cells.Item(180, 1).Value = 2
$ws.Cells.Item(181, 1).Value = 2
$ws.Cells.Item(182, 1).Value = 1
$ws.Cells.Item(183, 1).Value = 2
$ws.Cells.Item(184, 1).Value = 1
$ws.Cells.Item(185, 1).Value = 1
$ws.Cells.Item(186, 1).Value = 2
$ws.Cells.Item(187, 1).Value = 2
$ws.Cells.Item(189, 1).Value = 2
$ws.Cells.Item(192, 1).Value = 2
$ws.Cells.Item(193, 1).Value = 1
$ws.Cells.Item(195, 1).Value = 1
$ws.Cells.Item(196, 1).Value = 2
$ws.Cells.Item(197, 1).Value = 1
$ws.Cells.Item(201, 1).Value = 2
